# Generate Report for handback
#
# The "d189dbdd-312f-4618-a22f-7198e47b6c63.md" file has now been handed
# back (in sync with en-US) for both locales. Update the per-locale status
# and record the new handback timestamps, and reflect the updated status
# on the Overview sheet as well.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: row for d189dbdd-...md is row 3 (B = zh-cn, C = de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- zh-cn sheet: row 3 is the d189dbdd-...md entry ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $newStatus
$wsZhCn.Range("G3").Value = "2016-01-28 09:14:06"

# --- de-de sheet: row 3 is the d189dbdd-...md entry ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $newStatus
$wsDeDe.Range("G3").Value = "2016-01-28 09:14:26"
